$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reconnect Nets after length tuning: updated tuned lengths for LDQS+/LDQS-
$ws.Range("B22").Value = 28.91
$ws.Range("B23").Value = 28.91

# Leave the view scrolled/selected on the cell that was just edited
$ws.Range("B23").Select() | Out-Null
